$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) ---
# "Volume 31   Number  49" -> "...50"
$ws.Range("A8").Value = "Volume 31   Number  49" -replace "49$", "50"
# "Report Covering the Week  12/2/2024  Through  12/8/2024" -> new dates
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# Helper: convert a cell to "text zero" matching the style of C14/D15 (s=13, shared string "0")
function Set-TextZero($addr) {
    $ws.Range($addr).Value = "'0"
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Helper: convert a cell to the text "***.*" matching style s=13
function Set-TextStar($addr) {
    $ws.Range($addr).Value = "'***.*"
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Helper: set a numeric value while forcing the style to match a donor cell (paste formats only)
function Set-NumWithStyle($addr, $val, $donor) {
    $ws.Range($addr).Value = $val
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------- Row 15 ----------------
Set-TextZero("C15")
$ws.Range("F15").Value = 3
$ws.Range("M15").Value = 111.764705882353
$ws.Range("N15").Value = 12.5

# ---------------- Row 16 ----------------
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -57.142857142857
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 29
$ws.Range("H16").Value = -37.931034482758
$ws.Range("I16").Value = 284
$ws.Range("J16").Value = 316
$ws.Range("K16").Value = -10.126582278481
$ws.Range("L16").Value = -3.401360544217
$ws.Range("M16").Value = 39.215686274509
$ws.Range("N16").Value = -75.496117342536

# ---------------- Row 17 ----------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 16.129032258064
$ws.Range("I17").Value = 479
$ws.Range("J17").Value = 435
$ws.Range("K17").Value = 10.114942528735
$ws.Range("L17").Value = 45.592705167173
$ws.Range("M17").Value = 128.095238095238
$ws.Range("N17").Value = 40.882352941176

# ---------------- Row 18 ----------------
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 42.857142857142
$ws.Range("F18").Value = 33
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = -19.512195121951
$ws.Range("I18").Value = 414
$ws.Range("J18").Value = 534
$ws.Range("K18").Value = -22.471910112359
$ws.Range("L18").Value = -28.249566724436
$ws.Range("M18").Value = -10.775862068965
$ws.Range("N18").Value = -81.992170508916

# ---------------- Row 19 ----------------
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -47.826086956521
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 95
$ws.Range("H19").Value = -27.368421052631
$ws.Range("I19").Value = 1217
$ws.Range("J19").Value = 1302
$ws.Range("K19").Value = -6.528417818740
$ws.Range("L19").Value = -15.191637630662
$ws.Range("M19").Value = 85.801526717557
$ws.Range("N19").Value = -7.733131159969

# ---------------- Row 20 ----------------
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 36
$ws.Range("I20").Value = 512
$ws.Range("J20").Value = 489
$ws.Range("K20").Value = 4.703476482617
$ws.Range("L20").Value = 82.206405693950
$ws.Range("M20").Value = 82.206405693950
$ws.Range("N20").Value = -87.090267271810

# ---------------- Row 21 ----------------
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 59
$ws.Range("E21").Value = -22.033898305084
$ws.Range("F21").Value = 193
$ws.Range("G21").Value = 221
$ws.Range("H21").Value = -12.669683257918
$ws.Range("I21").Value = 2946
$ws.Range("J21").Value = 3112
$ws.Range("K21").Value = -5.334190231362
$ws.Range("L21").Value = -0.540175557056
$ws.Range("M21").Value = 60.282916213275
$ws.Range("N21").Value = -67.718606180144

# ---------------- Row 22 ----------------
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("J22").Value = 34
$ws.Range("K22").Value = -11.764705882352
$ws.Range("L22").Value = -25

# ---------------- Row 24 ----------------
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = 47.727272727272
$ws.Range("F24").Value = 236
$ws.Range("G24").Value = 216
$ws.Range("H24").Value = 9.259259259259
$ws.Range("I24").Value = 2691
$ws.Range("J24").Value = 2702
$ws.Range("K24").Value = -0.407105847520
$ws.Range("L24").Value = 0.635751682872
$ws.Range("M24").Value = 77.623762376237

# ---------------- Row 25 ----------------
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 33
$ws.Range("E25").Value = 21.212121212121
$ws.Range("F25").Value = 140
$ws.Range("G25").Value = 150
$ws.Range("H25").Value = -6.666666666666
$ws.Range("I25").Value = 1752
$ws.Range("J25").Value = 1512
$ws.Range("K25").Value = 15.873015873015
$ws.Range("L25").Value = 23.991507430997

# ---------------- Row 26 ----------------
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = 48.214285714285
$ws.Range("I26").Value = 947
$ws.Range("J26").Value = 852
$ws.Range("K26").Value = 11.150234741784
$ws.Range("L26").Value = 39.675516224188
$ws.Range("M26").Value = 37.845705967976

# ---------------- Row 27 ----------------
Set-TextZero("C27")
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200

# ---------------- Row 28 ----------------
Set-NumWithStyle "C28" 4 "F28"
Set-TextZero("D28")
Set-TextStar("E28")
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 98
$ws.Range("K28").Value = 4.255319148936
$ws.Range("L28").Value = 24.050632911392

# ---------------- Row 31 ----------------
Set-NumWithStyle "F31" 1 "I31"
$ws.Range("I31").Value = 11
$ws.Range("K31").Value = 37.5
$ws.Range("L31").Value = -15.384615384615
